$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (workbook view) ---
try { $wb.Windows.Item(1).Left = -33520 } catch {}
try { $wb.Windows.Item(1).Top = -8600 } catch {}

# --- Column widths ---
$ws.Columns("F").ColumnWidth = 121.333333
$ws.Columns("G").ColumnWidth = 71.5

# --- Row 17 height ---
$ws.Rows("17").RowHeight = 121

# --- D/E "X" marker cells (rows 17-22) ---
$ws.Range("E17").Value = "X"
$ws.Range("D18").Value = "X"
$ws.Range("E18").Value = "X"
$ws.Range("D19").Value = "X"
$ws.Range("E19").Value = "X"
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"
$ws.Range("D21").Value = "X"
$ws.Range("D22").Value = "X"

# --- New note / comment text values (order matters for shared-string indices) ---
$ws.Range("F17").Value = '<property> ::= [<visibility>] [‘/’] <name> [‘:’ <prop-type>] [‘[‘ <multiplicity-range> ’]’] [‘=’ <default>] [‘{‘ <prop-modifier > [‘,’ <prop-modifier >]* ’}’]'
$ws.Range("M18").Value = '<ownedComment xmi:type=''uml:Comment'' xmi:id=''_18_0_2_9cd0221_1415387037068_441145_10971'' body=''&lt;html&gt;&#10;  &lt;head&gt;&#10;&#9;&#9;&lt;style&gt;&#10;&#9;&#9;&#9;p {padding:0px; margin:0px;}&#10;&#9;&#9;&lt;/style&gt;&#10;&#9;&lt;/head&gt;&#10;  &lt;body&gt;&#10;    &lt;p&gt;&#10;The code system or code system or code system version that contained a description of the terminology code at the point in time it was referenced.&#10;&#10;    &lt;/p&gt;&#10;&#10;&lt;/body&gt;&#10;&lt;/html&gt;''>'
$ws.Range("N19").Value = '<annotatedElement xmi:idref=''_18_0_2_9cd0221_1415314303371_203085_10072''/>'
$ws.Range("M20").Value = '</ownedComment>'
$ws.Range("M21").Value = '<lowerValue xmi:type=''uml:LiteralInteger'' xmi:id=''_18_0_2_9cd0221_1415314345979_807009_10074''/>'
$ws.Range("M22").Value = '<xmi:Extension extender=''MagicDraw UML 18.0''>'
$ws.Range("N23").Value = '<modelExtension>'
$ws.Range("O24").Value = '<upperValue xmi:type=''uml:LiteralUnlimitedNatural'' xmi:id=''_18_0_2_9cd0221_1415314345979_723586_10075'' value=''1''/>'
$ws.Range("N25").Value = '</modelExtension>'
$ws.Range("M26").Value = '</xmi:Extension>'
$ws.Range("L27").Value = '</ownedAttribute>'

# --- G17: long block, wrap text style, set LAST so its shared-string index lands after the others ---
$ws.Range("G17").WrapText = $true
$ws.Range("G17").Value = '<ownedAttribute xmi:type=''uml:Property'' xmi:id=''_18_0_2_9cd0221_1415314303371_203085_10072'' name=''terminologyVersion'' visibility=''public'' type=''_18_0_2_9cd0221_1414186321748_769176_17132''>
      <ownedComment xmi:type=''uml:Comment'' xmi:id=''_18_0_2_9cd0221_1415387037068_441145_10971'' body=''&lt;html&gt;&#10;  &lt;head&gt;&#10;&#9;&#9;&lt;style&gt;&#10;&#9;&#9;&#9;p {padding:0px; margin:0px;}&#10;&#9;&#9;&lt;/style&gt;&#10;&#9;&lt;/head&gt;&#10;  &lt;body&gt;&#10;    &lt;p&gt;&#10;The code system or code system or code system version that contained a description of the terminology code at the point in time it was referenced.&#10;&#10;    &lt;/p&gt;&#10;&#10;&lt;/body&gt;&#10;&lt;/html&gt;''>
       <annotatedElement xmi:idref=''_18_0_2_9cd0221_1415314303371_203085_10072''/>
      </ownedComment>
      <lowerValue xmi:type=''uml:LiteralInteger'' xmi:id=''_18_0_2_9cd0221_1415314345979_807009_10074''/>
      <xmi:Extension extender=''MagicDraw UML 18.0''>
       <modelExtension>
        <upperValue xmi:type=''uml:LiteralUnlimitedNatural'' xmi:id=''_18_0_2_9cd0221_1415314345979_723586_10075'' value=''1''/>
       </modelExtension>
      </xmi:Extension>
     </ownedAttribute>'

# --- Selection ---
$ws.Range("F20").Select() | Out-Null

